# Auto-generated edit script: update crafting-profit market-price data cells
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(128, 8).Value = 17446.5
$ws.Cells.Item(128, 10).Value = 17446.5
$ws.Cells.Item(128, 12).Value = 17446.5
$ws.Cells.Item(128, 14).Value = -27406.5
$ws.Cells.Item(132, 8).Value = 1049.5714
$ws.Cells.Item(132, 9).Value = 609.6429000000001
$ws.Cells.Item(132, 10).Value = 1929.4286
$ws.Cells.Item(132, 11).Value = 1828.9287
$ws.Cells.Item(132, 12).Value = 5788.2858
$ws.Cells.Item(132, 13).Value = 701.0712999999998
$ws.Cells.Item(132, 14).Value = -10848.2858
$ws.Cells.Item(141, 8).Value = 3760
$ws.Cells.Item(141, 9).Value = 2796
$ws.Cells.Item(141, 10).Value = 5366.6665
$ws.Cells.Item(141, 11).Value = 8388
$ws.Cells.Item(141, 12).Value = 16099.9995
$ws.Cells.Item(141, 13).Value = -3208
$ws.Cells.Item(141, 14).Value = -26459.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3690.1714
$ws.Cells.Item(61, 9).Value = 2726.5715
$ws.Cells.Item(61, 10).Value = 5135.5713
$ws.Cells.Item(61, 11).Value = 2726.5715
$ws.Cells.Item(61, 12).Value = 5135.5713
$ws.Cells.Item(61, 13).Value = -2514.5715
$ws.Cells.Item(61, 14).Value = -5559.5713
$ws.Cells.Item(74, 8).Value = 1284.2812
$ws.Cells.Item(74, 9).Value = 909.381
$ws.Cells.Item(74, 10).Value = 2000
$ws.Cells.Item(74, 11).Value = 909.381
$ws.Cells.Item(74, 12).Value = 2000
$ws.Cells.Item(74, 13).Value = -35.38099999999997
$ws.Cells.Item(74, 14).Value = -3748
$ws.Cells.Item(77, 8).Value = 1284.2812
$ws.Cells.Item(77, 9).Value = 909.381
$ws.Cells.Item(77, 10).Value = 2000
$ws.Cells.Item(77, 11).Value = 4546.905
$ws.Cells.Item(77, 12).Value = 10000
$ws.Cells.Item(77, 13).Value = -178.9049999999997
$ws.Cells.Item(77, 14).Value = -18736
$ws.Cells.Item(122, 8).Value = 2012.0476
$ws.Cells.Item(122, 9).Value = 1890.875
$ws.Cells.Item(122, 10).Value = 2399.8
$ws.Cells.Item(122, 11).Value = 5672.625
$ws.Cells.Item(122, 12).Value = 7199.400000000001
$ws.Cells.Item(122, 13).Value = -3222.625
$ws.Cells.Item(122, 14).Value = -12099.4
$ws.Cells.Item(128, 8).Value = 41750
$ws.Cells.Item(128, 10).Value = 41750
$ws.Cells.Item(128, 12).Value = 41750
$ws.Cells.Item(128, 14).Value = -51710
$ws.Cells.Item(136, 8).Value = 3690.1714
$ws.Cells.Item(136, 9).Value = 2726.5715
$ws.Cells.Item(136, 10).Value = 5135.5713
$ws.Cells.Item(136, 11).Value = 8179.7145
$ws.Cells.Item(136, 12).Value = 15406.7139
$ws.Cells.Item(136, 13).Value = -5629.7145
$ws.Cells.Item(136, 14).Value = -20506.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 750.5
$ws.Cells.Item(107, 9).Value = 750.5
$ws.Cells.Item(107, 11).Value = 750.5
$ws.Cells.Item(107, 13).Value = 1169.5
$ws.Cells.Item(131, 8).Value = 48923.332
$ws.Cells.Item(131, 10).Value = 48923.332
$ws.Cells.Item(131, 12).Value = 48923.332
$ws.Cells.Item(131, 14).Value = -59003.332
$ws.Cells.Item(134, 8).Value = 2831.0667
$ws.Cells.Item(134, 9).Value = 1524.7273
$ws.Cells.Item(134, 10).Value = 6423.5
$ws.Cells.Item(134, 11).Value = 4574.1819
$ws.Cells.Item(134, 12).Value = 19270.5
$ws.Cells.Item(134, 13).Value = -2039.1819
$ws.Cells.Item(134, 14).Value = -24340.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1590.4286
$ws.Cells.Item(16, 9).Value = 1605.5
$ws.Cells.Item(16, 10).Value = 1500
$ws.Cells.Item(16, 11).Value = 1605.5
$ws.Cells.Item(16, 12).Value = 1500
$ws.Cells.Item(16, 13).Value = -1318.5
$ws.Cells.Item(16, 14).Value = -2074
$ws.Cells.Item(31, 8).Value = 4140.0957
$ws.Cells.Item(31, 9).Value = 825.7241
$ws.Cells.Item(31, 10).Value = 6324.5684
$ws.Cells.Item(31, 11).Value = 825.7241
$ws.Cells.Item(31, 12).Value = 6324.5684
$ws.Cells.Item(31, 13).Value = -530.7241
$ws.Cells.Item(31, 14).Value = -6914.5684
$ws.Cells.Item(34, 8).Value = 4140.0957
$ws.Cells.Item(34, 9).Value = 825.7241
$ws.Cells.Item(34, 10).Value = 6324.5684
$ws.Cells.Item(34, 11).Value = 825.7241
$ws.Cells.Item(34, 12).Value = 6324.5684
$ws.Cells.Item(34, 13).Value = -623.7241
$ws.Cells.Item(34, 14).Value = -6728.5684
$ws.Cells.Item(58, 8).Value = 668.5238000000001
$ws.Cells.Item(58, 9).Value = 682.45
$ws.Cells.Item(58, 10).Value = 390
$ws.Cells.Item(58, 11).Value = 682.45
$ws.Cells.Item(58, 12).Value = 390
$ws.Cells.Item(58, 13).Value = -479.45
$ws.Cells.Item(58, 14).Value = -796
$ws.Cells.Item(113, 8).Value = 1590.4286
$ws.Cells.Item(113, 9).Value = 1605.5
$ws.Cells.Item(113, 10).Value = 1500
$ws.Cells.Item(113, 11).Value = 1605.5
$ws.Cells.Item(113, 12).Value = 1500
$ws.Cells.Item(113, 13).Value = 564.5
$ws.Cells.Item(113, 14).Value = -5840
$ws.Cells.Item(122, 8).Value = 1081.8182
$ws.Cells.Item(122, 9).Value = 1081.8182
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 3245.4546
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -795.4546
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 2834.682
$ws.Cells.Item(132, 9).Value = 2515.9285
$ws.Cells.Item(132, 10).Value = 3392.5
$ws.Cells.Item(132, 11).Value = 7547.7855
$ws.Cells.Item(132, 12).Value = 10177.5
$ws.Cells.Item(132, 13).Value = -5017.7855
$ws.Cells.Item(132, 14).Value = -15237.5
$ws.Cells.Item(134, 8).Value = 4198.5454
$ws.Cells.Item(134, 9).Value = 5729.5264
$ws.Cells.Item(134, 10).Value = 2120.7856
$ws.Cells.Item(134, 11).Value = 17188.5792
$ws.Cells.Item(134, 12).Value = 6362.3568
$ws.Cells.Item(134, 13).Value = -14653.5792
$ws.Cells.Item(134, 14).Value = -11432.3568
$ws.Cells.Item(136, 8).Value = 668.5238000000001
$ws.Cells.Item(136, 9).Value = 682.45
$ws.Cells.Item(136, 10).Value = 390
$ws.Cells.Item(136, 11).Value = 2047.35
$ws.Cells.Item(136, 12).Value = 1170
$ws.Cells.Item(136, 13).Value = 502.6499999999999
$ws.Cells.Item(136, 14).Value = -6270

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 753.10205
$ws.Cells.Item(107, 9).Value = 808.6087
$ws.Cells.Item(107, 10).Value = 704
$ws.Cells.Item(107, 11).Value = 2425.8261
$ws.Cells.Item(107, 12).Value = 2112
$ws.Cells.Item(107, 13).Value = -505.8261000000002
$ws.Cells.Item(107, 14).Value = -5952
$ws.Cells.Item(131, 8).Value = 1460.42
$ws.Cells.Item(131, 9).Value = 781.6667
$ws.Cells.Item(131, 10).Value = 1503.7446
$ws.Cells.Item(131, 11).Value = 2345.0001
$ws.Cells.Item(131, 12).Value = 4511.2338
$ws.Cells.Item(131, 13).Value = 2694.9999
$ws.Cells.Item(131, 14).Value = -14591.2338
$ws.Cells.Item(136, 8).Value = 1526.6923
$ws.Cells.Item(136, 9).Value = 1487.25
$ws.Cells.Item(136, 10).Value = 2000
$ws.Cells.Item(136, 11).Value = 4461.75
$ws.Cells.Item(136, 12).Value = 6000
$ws.Cells.Item(136, 13).Value = 638.25
$ws.Cells.Item(136, 14).Value = -16200
$ws.Cells.Item(137, 8).Value = 1558.8077
$ws.Cells.Item(137, 9).Value = 1292.8572
$ws.Cells.Item(137, 10).Value = 1869.0834
$ws.Cells.Item(137, 11).Value = 3878.5716
$ws.Cells.Item(137, 12).Value = 5607.2502
$ws.Cells.Item(137, 13).Value = 1221.4284
$ws.Cells.Item(137, 14).Value = -15807.2502

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 2241.6667
$ws.Cells.Item(113, 9).Value = 2275
$ws.Cells.Item(113, 10).Value = 2175
$ws.Cells.Item(113, 11).Value = 2275
$ws.Cells.Item(113, 12).Value = 2175
$ws.Cells.Item(113, 13).Value = -105
$ws.Cells.Item(113, 14).Value = -6515
$ws.Cells.Item(122, 8).Value = 2336.6487
$ws.Cells.Item(122, 9).Value = 1898.4814
$ws.Cells.Item(122, 10).Value = 3519.7
$ws.Cells.Item(122, 11).Value = 5695.4442
$ws.Cells.Item(122, 12).Value = 10559.1
$ws.Cells.Item(122, 13).Value = -3245.4442
$ws.Cells.Item(122, 14).Value = -15459.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 5467.0356
$ws.Cells.Item(132, 9).Value = 5496.263
$ws.Cells.Item(132, 10).Value = 5405.3335
$ws.Cells.Item(132, 11).Value = 16488.789
$ws.Cells.Item(132, 12).Value = 16216.0005
$ws.Cells.Item(132, 13).Value = -13958.789
$ws.Cells.Item(132, 14).Value = -21276.0005
$ws.Cells.Item(136, 8).Value = 8132173
$ws.Cells.Item(136, 9).Value = 781.1111
$ws.Cells.Item(136, 10).Value = 14495871
$ws.Cells.Item(136, 11).Value = 2343.3333
$ws.Cells.Item(136, 12).Value = 43487613
$ws.Cells.Item(136, 13).Value = 206.6667000000002
$ws.Cells.Item(136, 14).Value = -43492713

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(110, 8).Value = 51577.5
$ws.Cells.Item(110, 10).Value = 51577.5
$ws.Cells.Item(110, 12).Value = 51577.5
$ws.Cells.Item(110, 14).Value = -59757.5
$ws.Cells.Item(122, 8).Value = 51468.3
$ws.Cells.Item(122, 9).Value = 63741.312
$ws.Cells.Item(122, 10).Value = 2376.25
$ws.Cells.Item(122, 11).Value = 191223.936
$ws.Cells.Item(122, 12).Value = 7128.75
$ws.Cells.Item(122, 13).Value = -188773.936
$ws.Cells.Item(122, 14).Value = -12028.75
$ws.Cells.Item(136, 8).Value = 4802.7144
$ws.Cells.Item(136, 9).Value = 847.625
$ws.Cells.Item(136, 10).Value = 5974.593
$ws.Cells.Item(136, 11).Value = 2542.875
$ws.Cells.Item(136, 12).Value = 17923.779
$ws.Cells.Item(136, 13).Value = 7.125
$ws.Cells.Item(136, 14).Value = -23023.779

